$d = $word.ActiveDocument

# Change 1: paragraph about recommended usage guidelines -> Italian translation
$rng1 = $d.Content
$rng1.Find.Execute(
    "These are recommended usage guidelines for maintaining a consistent design aesthetic for the SmartCash brand." + [char]160 + "A strong and consistent visual identity of our logo will help keep a consistent look, recognition and familiarity now and in the future. Standardization of colours will go a long way to enforce a reliable and positive impression to our identity in the blockchain space.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0
)
$rng1.Text = "Queste sono le linee guida di utilizzo consigliate per mantenere un'estetica di design coerente per il marchio SmartCash. Un'identità visiva forte e coerente del nostro logo contribuirà a mantenere un aspetto coerente, riconoscibile e familiare per il presente e per il futuro. La standardizzazione dei colori sarà d'aiuto per rafforzare positivamente l'affidabilità della nostra identità nello spazio della blockchain."

# Change 2: "Official font is" + nbsp -> "Il carattere ufficiale è " (regular space)
$rng2 = $d.Content
$rng2.Find.Execute(
    "Official font is" + [char]160,
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0
)
$rng2.Text = "Il carattere ufficiale è "

# Change 3: "Source Sans Pro" -> " Source Sans Pro " (inside hyperlink run)
$rng3 = $d.Content
$rng3.Find.Execute(
    "Source Sans Pro",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0
)
$origStart3 = $rng3.Start
$origEnd3 = $rng3.End
$newText3 = " Source Sans Pro "
$rng3.Text = $newText3

# Restore underline formatting on exactly the new text (a plain re-Find here
# can mis-extend into the following run, so address the range by its known
# character positions instead).
$newEnd3 = $origStart3 + $newText3.Length
$fix3 = $d.Range($origStart3, $newEnd3)
$fix3.Font.Underline = 1

# Change 4: nbsp + "/ PT Sans" -> regular space + "/ PT Sans"
$rng4 = $d.Content
$rng4.Find.Execute(
    [char]160 + "/ PT Sans",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0
)
$rng4.Text = " / PT Sans"
